# edit.ps1
# Adds a new survey wave column (11. 11. 2021) to both worksheets of the
# "ZBP_03_strategie_domacnosti" workbook ("data" and "pocetR"), and
# refreshes the trailing "aktualizace" (updated-on) note on each sheet
# from 20. 10. 2021 to 18. 11. 2021.

$wb = $excel.ActiveWorkbook

# xlPasteFormats - used below to clone a header cell's style onto the
# freshly added header cell without hand-rebuilding font/border/align.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheet "data" — percentages. New column AJ (column 36), header date
# "11. 11. 2021", formatted like the rest of row 1 (bordered, bold,
# centered header style already used by AI1).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

$ws1.Range("AI1").Copy()
$ws1.Range("AJ1").PasteSpecial($xlPasteFormats)
$ws1.Range("AJ1").Value = "11. 11. 2021"

$data_AJ = 0.22,0.12,0.46,0.29,0.17,0.22,0.24,0.23,0.2,0.21,0.23,0.45,0.2,0.22,0.23,0.21,0.2,0.28,0.2,0.14,0.16,0.23,0.44,0.39,0.11,0.09,0.12,0.19,0.11,0.1,0.12,0.18,0.18,0.1,0.13,0.13,0.07,0.18,0.14,0.11,0.07,0.06,0.19,0.27
for ($i = 0; $i -lt $data_AJ.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 36).Value = $data_AJ[$i]
}

$ws1.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 18. 11. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR" — sample sizes. New column AI (column 35), header date
# "11. 11. 2021", formatted like the rest of row 1 (style used by AH1).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

$ws2.Range("AH1").Copy()
$ws2.Range("AI1").PasteSpecial($xlPasteFormats)
$ws2.Range("AI1").Value = "11. 11. 2021"

$pocet_AI = 1709,167,365,1177,820,152,489,248,783,140,103,683,788,590,331,182,618,536,250,506,294,155
for ($i = 0; $i -lt $pocet_AI.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 35).Value = $pocet_AI[$i]
}

# Row 24 is a blank formatting row under the table; keep column AI
# consistent with it (empty cell), same as the rest of that row. Simply
# assigning "" is a no-op for a brand-new cell, so touch a harmless
# no-op border property instead to materialize the (still contentless,
# unstyled) cell in the sheet.
$ws2.Range("AI24").Borders.LineStyle = -4142

$ws2.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 18. 11. 2021"
